$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "VALOR MORA" total: 44640 -> 4640
$ws.Range("E11").Value = 4640

# Cant. Trabajadores / Cant. Periodos: 2 -> 1 (the KENIS PAYARES worker entry is being removed)
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Salario Basico for the remaining worker (DUGLAN JOSE HERNANDEZ MATOS) updated
$ws.Range("G16").Value = 1423500

# Remove the second worker's row (CC 1047393798 - KENIS JUVENAL PAYARES BARRIOS - 2506)
# entirely, shifting the signature block (rows 22-23) up to rows 21-22.
$ws.Rows.Item(17).Delete()
